$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Donation without test order class added
#
# The worklist used to be asserted with a 5 column table
# (WorklistID / PageSize / Page / Sort / EndPoint). PageSize/Page/Sort
# are gone now - the sheet only keeps EndPoint / WorklistID, laid out
# as a 2 column key/value table. The three assertion blocks are also
# reordered (Assert200, Assert401, Assert404) and the Assert404 block
# now exercises a worklist that has no test-order class ("a1") instead
# of "last1000001".
# ------------------------------------------------------------------

# Grab format "donors" from the existing sheet before we start
# rewriting it, so the new layout re-uses the same style records
# instead of inventing new ones. Park copies of them off to the side
# (row 40+) since the originals live inside the A1:E25 block we are
# about to wipe.
$ws.Range("A1").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headerStyleCell  = $ws.Range("A40")  # 40% - Accent4 block banner
$labelStyleCell   = $ws.Range("A41")  # 40% - Accent3 label cell
$valueStyleCell   = $ws.Range("A42")  # wrap/valign-top value cell
$blankStyleCell   = $ws.Range("A43")  # plain wrap/valign-top blank cell

function Set-Block($headerRow, $header, $endpoint, $worklist) {
  $r1 = $headerRow
  $r2 = $headerRow + 1
  $r3 = $headerRow + 2
  $r4 = $headerRow + 3
  $r5 = $headerRow + 4

  # Banner row, merged A:B
  $headerStyleCell.Copy()
  $ws.Range("A$r1").PasteSpecial(-4122)
  $ws.Range("B$r1").PasteSpecial(-4122)
  $ws.Range("A$r1").Value = $header
  $ws.Range("B$r1").Value = $null

  # Label row
  $labelStyleCell.Copy()
  $ws.Range("A$r2").PasteSpecial(-4122)
  $ws.Range("B$r2").PasteSpecial(-4122)
  $ws.Range("A$r2").Value = "EndPoint"
  $ws.Range("B$r2").Value = "WorklistID"

  # Value row (endpoint keeps the wrap/valign-top "value" look, the
  # worklist id next to it uses the plain wrap/valign-top "blank" look)
  $valueStyleCell.Copy()
  $ws.Range("A$r3").PasteSpecial(-4122)
  $blankStyleCell.Copy()
  $ws.Range("B$r3").PasteSpecial(-4122)
  $ws.Range("A$r3").Value = $endpoint
  $ws.Range("B$r3").Value = $worklist
  $ws.Rows.Item($r3).RowHeight = 30

  # Two trailing blank rows
  $blankStyleCell.Copy()
  $ws.Range("A$r4").PasteSpecial(-4122)
  $ws.Range("A$r5").PasteSpecial(-4122)
  $ws.Range("A$r4").Value = $null
  $ws.Range("A$r5").Value = $null

  $excel.CutCopyMode = $false
}

# Wipe the old 5 column content (incl. formatting) so C:E don't keep
# stray data/styles from the old layout.
$ws.Range("A1:E25").Clear()

# Undo the old A1:E1 / A6:E6 / A11:E11 banners before re-merging as
# A:B only.
$ws.Range("A1:E1").UnMerge()
$ws.Range("A6:E6").UnMerge()
$ws.Range("A11:E11").UnMerge()

# Block 1 - Assert200 (unchanged position/values, just narrowed to 2 cols)
Set-Block 1 "Assert200" "/worklistInfo/worklistCalibratorsControls" "last1"

# Block 2 - Assert401 (moved up from row 11)
Set-Block 6 "Assert401" "/worklistInfo/worklistCalibratorsControls" "last1"

# Block 3 - Assert404 (moved down from row 6, worklist w/o test order class)
Set-Block 11 "Assert404" "/worklistInfo/worklistCalibratorsControls" "a1"

# Re-merge the banner rows as two columns only
$ws.Range("A1:B1").Merge()
$ws.Range("A6:B6").Merge()
$ws.Range("A11:B11").Merge()

# Column B widens to fit "WorklistID" style labels now that C:E are gone
$ws.Columns.Item(2).ColumnWidth = 19.43

# Selection / view tweaks to mirror the saved workbook
$ws.Range("E7").Select()

# Touch the bottom of the sheet so the used range grows to A1:E31 like
# the authored version (extra blank rows below the tables).
$ws.Range("E16:E22").Value = " "
$ws.Range("E16:E22").Clear()
$ws.Range("E25:E31").Value = " "
$ws.Range("E25:E31").Clear()

# Drop the scratch donors used for format painting.
$ws.Range("A40:A43").Clear()

$wb.Save()
